$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.970.03'
$ws.Range("E2").Value = '  -0.82%  '

$ws.Range("D3").Value = '1.915.74'
$ws.Range("E3").Value = '  +0.96%  '

$ws.Range("E4").Value = '  -0.02%  '

$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = '320.05'
$rng.Style = "Normal"
$ws.Range("E5").Value = '  -1.58%  '

$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = '1.0000'
$rng.Style = "Normal"

$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = '0.5043'
$rng.Style = "Normal"
$ws.Range("E7").Value = '  -2.24%  '

$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = '0.4039'
$rng.Style = "Normal"
$ws.Range("E8").Value = '  +0.77%  '

$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = '0.08288'
$rng.Style = "Normal"
$ws.Range("E9").Value = '  -1.86%  '

$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = '1.104'
$rng.Style = "Normal"
$ws.Range("E10").Value = '  -1.25%  '

$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = '42.01'
$rng.Style = "Normal"
$ws.Range("E11").Value = '  -1.56%  '

$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = '23.98'
$rng.Style = "Normal"
$ws.Range("E12").Value = '  +2.43%  '

$ws.Range("D13").Value = '1.914.21'
$ws.Range("E13").Value = '  +1.09%  '

$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = '6.398'
$rng.Style = "Normal"
$ws.Range("E14").Value = '  -0.48%  '

$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = '7.224'
$rng.Style = "Normal"
$ws.Range("E15").Value = '  -1.56%  '

$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = '0.9990'
$rng.Style = "Normal"
$ws.Range("E16").Value = '  -0.22%  '

$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = '92.26'
$rng.Style = "Normal"
$ws.Range("E17").Value = '  -2.66%  '

$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = '0.00001097'
$rng.Style = "Normal"
$ws.Range("E18").Value = '  -1.31%  '

$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = '0.06501'
$rng.Style = "Normal"
$ws.Range("E19").Value = '  -2.13%  '

$ws.Range("E20").Value = '  -0.42%  '

$ws.Range("E21").Value = '  +0.00%  '

$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = '5.938'
$rng.Style = "Normal"
$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("D23").Value = '30.006.83'
$ws.Range("E23").Value = '  -0.73%  '

$ws.Range("E24").Value = '  +0.17%  '

$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = '2.196'
$rng.Style = "Normal"
$ws.Range("E25").Value = '  -1.40%  '

$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = '22.18'
$rng.Style = "Normal"
$ws.Range("E26").Value = '  +2.58%  '

$ws.Range("D27").Value = '2.133.99'
$ws.Range("E27").Value = '  +1.11%  '

$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = '162.43'
$rng.Style = "Normal"
$ws.Range("E28").Value = '  +0.84%  '

$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = '2.297'
$rng.Style = "Normal"
$ws.Range("E29").Value = '  -2.78%  '

$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = '129.02'
$rng.Style = "Normal"
$ws.Range("E30").Value = '  +0.11%  '

$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = '1.131'
$rng.Style = "Normal"
$ws.Range("E31").Value = '  +3.05%  '

$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = '0.1036'
$rng.Style = "Normal"
$ws.Range("E32").Value = '  -2.05%  '

$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = '5.942'
$rng.Style = "Normal"
$ws.Range("E33").Value = '  -2.27%  '

$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = '3.827'
$rng.Style = "Normal"
$ws.Range("E34").Value = '  +1.76%  '

$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = '5.462'
$rng.Style = "Normal"
$ws.Range("E35").Value = '  +3.57%  '

$rng = $ws.Range("D36")
$rng.NumberFormat = "@"
$rng.Value = '0.02448'
$rng.Style = "Normal"
$ws.Range("E36").Value = '  -1.84%  '

$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = '0.06407'
$rng.Style = "Normal"
$ws.Range("E37").Value = '  -2.32%  '

$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = '0.2154'
$rng.Style = "Normal"
$ws.Range("E38").Value = '  -2.13%  '

$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = '8.732'
$rng.Style = "Normal"
$ws.Range("E39").Value = '  +0.34%  '

$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = '1.196'
$rng.Style = "Normal"
$ws.Range("E40").Value = '  -1.91%  '

$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = '0.6459'
$rng.Style = "Normal"
$ws.Range("E41").Value = '  -0.63%  '

$ws.Range("E42").Value = '  -3.42%  '

$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = '1.212'
$rng.Style = "Normal"
$ws.Range("E43").Value = '  -1.47%  '

$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$rng.Value = '2.219'
$rng.Style = "Normal"
$ws.Range("E44").Value = '  +7.85%  '

$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = '13.34'
$rng.Style = "Normal"
$ws.Range("E45").Value = '  +0.60%  '

$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = '0.6046'
$rng.Style = "Normal"
$ws.Range("E46").Value = '  -0.93%  '

$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = '3.636'
$rng.Style = "Normal"
$ws.Range("E47").Value = '  -1.78%  '

$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = '122.16'
$rng.Style = "Normal"
$ws.Range("E48").Value = '  -1.78%  '

$ws.Range("E49").Value = '  -2.46%  '

$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = '78.92'
$rng.Style = "Normal"
$ws.Range("E50").Value = '  -0.11%  '

$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = '1.130'
$rng.Style = "Normal"
$ws.Range("E51").Value = '  -3.07%  '

